# AFNT Progress Calculator - spelling fix: "Partically" -> "Partially"
# (applies to both the "Partially Complete" column header and the
# "Total Partially Completed" summary row label), plus re-selecting the
# Non-Functional Requirements status rows on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header above the Functional Requirements table (row 1) and the
# Non-Functional Requirements table (row 8) both read "Partically
# Complete" - fix the typo in both places.
$ws.Range("D1").Value = "Partially Complete"
$ws.Range("D8").Value = "Partially Complete"

# Summary label in the totals block has the same typo.
$ws.Range("A23").Value = "Total Partially Completed"

# Update the active selection left on the sheet.
[void]$ws.Range("A11:A17").Select()
